$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh: (cell, new text value)
$updates = @(
    @("D2", "305.07"),
    @("E2", "1.33%"),
    @("D3", "35.91"),
    @("E3", "2.37%"),
    @("E4", "1.13%"),
    @("D5", "0.08084"),
    @("E5", "1.28%"),
    @("D6", "1.936"),
    @("E6", "0.72%"),
    @("D7", "4.180"),
    @("E7", "3.45%"),
    @("D8", "7.758"),
    @("E8", "0.28%"),
    @("D9", "0.9263"),
    @("E9", "0.27%"),
    @("D10", "0.1350"),
    @("E10", "2.41%"),
    @("D11", "0.1903"),
    @("E11", "2.60%"),
    @("D12", "0.09167"),
    @("E12", "-5.40%"),
    @("D13", "0.03408"),
    @("E13", "-5.74%"),
    @("D14", "0.09825"),
    @("E14", "-0.34%"),
    @("D15", "0.001404"),
    @("E15", "0.51%"),
    @("D16", "0.005905"),
    @("E16", "1.77%"),
    @("D17", "3.554"),
    @("E17", "1.50%"),
    @("E19", "0.60%"),
    @("E20", "1.67%"),
    @("D21", "4.900"),
    @("E21", "-3.20%"),
    @("D22", "0.2601"),
    @("E22", "5.45%"),
    @("D23", "0.04405"),
    @("E23", "-2.70%"),
    @("D24", "0.001221"),
    @("E24", "0.22%"),
    @("D25", "0.004807"),
    @("E25", "-0.02%"),
    @("E26", "3.95%"),
    @("D27", "0.0003131"),
    @("E27", "4.19%"),
    @("E39", "5.84%"),
    @("D40", "0.04916"),
    @("E40", "4.45%"),
    @("D41", "0.007621"),
    @("E41", "1.03%"),
    @("D42", "0.01023"),
    @("E42", "6.28%"),
    @("D43", "0.1374"),
    @("E43", "3.47%"),
    @("E44", "-0.52%"),
    @("E45", "0.41%"),
    @("D46", "0.00006386"),
    @("E46", "2.85%"),
    @("E47", "0.00%"),
    @("D48", "63.57"),
    @("E48", "-1.41%"),
    @("E49", "-20.00%"),
    @("D50", "0.00002102"),
    @("E50", "0.00%"),
    @("D51", "0.0002002"),
    @("E51", "0.00%")
)

foreach ($u in $updates) {
    $cell = $ws.Range($u[0])
    # Force text storage so values like "305.07" / "1.33%" stay strings,
    # matching the sheet's original inline-string cell type (not auto-coerced to numbers).
    $cell.NumberFormat = "@"
    $cell.Value = $u[1]
    $cell.Style = "Normal"
}
